$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (first worksheet)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1222
$ws1.Range("F3").Value = 14737
$ws1.Range("F4").Value = 18052
$ws1.Range("F5").Value = 18052
$ws1.Range("F7").Value = 87
$ws1.Range("F16").Value = 65
$ws1.Range("F17").Value = 177
$ws1.Range("F19").Value = 1359
$ws1.Range("F21").Value = 81
$ws1.Range("F22").Value = 73
$ws1.Range("F23").Value = 219
$ws1.Range("F24").Value = 7464
$ws1.Range("F25").Value = 985
$ws1.Range("F28").Value = 1196
$ws1.Range("F30").Value = 5895
$ws1.Range("F31").Value = 77
$ws1.Range("F32").Value = 50
$ws1.Range("F33").Value = 152
$ws1.Range("F35").Value = 242
$ws1.Range("F36").Value = 5182
$ws1.Range("F38").Value = 36

# Sheet 4: 全部类型 (fourth worksheet)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1222
$ws4.Range("F3").Value = 14737
$ws4.Range("F4").Value = 18052
$ws4.Range("F5").Value = 18052
$ws4.Range("F7").Value = 87
$ws4.Range("F16").Value = 65
$ws4.Range("F17").Value = 178
$ws4.Range("F19").Value = 1359
$ws4.Range("F21").Value = 81
$ws4.Range("F23").Value = 73
$ws4.Range("F24").Value = 219
$ws4.Range("F25").Value = 7464
$ws4.Range("F26").Value = 985
$ws4.Range("F29").Value = 1196
$ws4.Range("F32").Value = 5895
$ws4.Range("F33").Value = 77
$ws4.Range("F34").Value = 50
$ws4.Range("F35").Value = 152
$ws4.Range("F37").Value = 242
$ws4.Range("F38").Value = 5182
$ws4.Range("F40").Value = 36
